# Elegoo Material Properties.xlsx - add "PC" column (O) + footnote row
# Mirrors the author's edit: a new PC material column is appended after
# column N (TPU-95A/ASA/PAHT-CF... last col), plus a changelog note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Bring column O into the table: copy column N's layout (widths
#        come from the <cols> rule, but per-cell formatting/merges/number
#        formats need to be seeded) down to row 20 -------------------------
$ws.Range("N1:N20").Copy($ws.Range("O1:O20"))

# Column O should be the same fixed width as the other data columns (F:N).
$ws.Range("O1").EntireColumn.ColumnWidth = 16.6328125

# --- 2. Header (row 1/2), merged like every other material column --------
$ws.Range("O1:O2").Merge()
$ws.Range("O1").Value2 = "PC"

# --- 3. Extend the three section-title merged bars one column to the right
$ws.Range("A3:N3").UnMerge()
$ws.Range("A3:O3").Merge()

$ws.Range("A10:N10").UnMerge()
$ws.Range("A10:O10").Merge()

$ws.Range("A15:N15").UnMerge()
$ws.Range("A15:O15").Merge()

# --- 4. Print-settings block (rows 4-9) for the new PC column ------------
$ws.Range("O4").Value2  = "80 ± 5"          # Drying Temp      (same as N4)
$ws.Range("O5").Value2  = "<50 <20%"        # Temp & Humidity
$ws.Range("O6").Value2  = "260-290"         # Nozzle Temp      (same as N6)
$ws.Range("O7").Value2  = "90-110"          # Bed Temp
$ws.Range("O8").Value2  = "Textured or Other" # Plate Type
$ws.Range("O9").Value2  = "<100"            # Print Speed      (same as N9)

# Row 5 ("Temp & Humidity") no longer needs the tall autofit height that
# the old 5-column wrap required.
$ws.Rows.Item(5).RowHeight = 14

# --- 5. Physical properties block (rows 11-14) ----------------------------
$ws.Range("O11").Value2 = 1.2    # Density
$ws.Range("O12").Value2 = 228    # Melt Temp
$ws.Range("O13").Value2 = 119    # VICAT Soft Temp
$ws.Range("O14").Value2 = 109    # Heat Deflection Temp

# --- 6. Mechanical properties block (rows 16-20) --------------------------
$ws.Range("O16").Value2 = "56 ± 2"        # Tensile Strength
$ws.Range("O17").Value2 = "8.8% ± 1.9%"   # Breaking Elongation XY
$ws.Range("O18").Value2 = "2728 ± 113"    # Bending Modulus XY
$ws.Range("O19").Value2 = "114 ± 5"       # Bending Strength
$ws.Range("O20").Value2 = "797 ± 3.2"     # Impact Strength

# --- 7. Changelog footnote -------------------------------------------------
$ws.Range("A22").Value2 = "vs 2 - added PC - 18 Jul 25"

# --- 8. Print setup tweaks that came along with the re-layout ------------
$ws.PageSetup.Zoom = 58
$ws.PageSetup.CenterVertically = $false

# --- 9. Leave the selection where the author left it ----------------------
$ws.Range("L42").Select()
